$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 57.35848733333334
$ws.Range("N2").Value = 172.075462
$ws.Range("O2").Value = 0.261658309594631
$ws.Range("P2").Value = 0.261658309594631
$ws.Range("Q2").Value = 1302.075576426794
$ws.Range("R2").Value = 11718.68018784115
$ws.Range("S2").Value = 0.02254256690435615
$ws.Range("T2").Value = 0.02254256690435615
$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.2957894889638607
$ws.Range("P3").Value = 0.2957894889638607
$ws.Range("Q3").Value = 1471.920650791778
$ws.Range("R3").Value = 13247.285857126
$ws.Range("S3").Value = 0.02548305977709321
$ws.Range("T3").Value = 0.02548305977709321
$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 29.294891
$ws.Range("N4").Value = 87.88467299999999
$ws.Range("O4").Value = 0.1336376186888105
$ws.Range("P4").Value = 0.1336376186888105
$ws.Range("Q4").Value = 665.0133896229509
$ws.Range("R4").Value = 5985.120506606558
$ws.Range("S4").Value = 0.01151324016767691
$ws.Range("T4").Value = 0.01151324016767691
$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("M5").Value = 67.71760166666667
$ws.Range("N5").Value = 203.152805
$ws.Range("O5").Value = 0.3089145827526977
$ws.Range("P5").Value = 0.3089145827526977
$ws.Range("Q5").Value = 1537.234319168035
$ws.Range("R5").Value = 13835.10887251231
$ws.Range("S5").Value = 0.02661382189704722
$ws.Range("T5").Value = 0.02661382189704722
$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 57.35848733333334
$ws.Range("N6").Value = 172.075462
$ws.Range("O6").Value = 0.261658309594631
$ws.Range("P6").Value = 0.261658309594631
$ws.Range("Q6").Value = 7602.567033157186
$ws.Range("R6").Value = 68423.10329841467
$ws.Range("S6").Value = 0.1316216808705603
$ws.Range("T6").Value = 0.1316216808705603
$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.2957894889638607
$ws.Range("P7").Value = 0.2957894889638607
$ws.Range("S7").Value = 0.1487906490781145
$ws.Range("T7").Value = 0.1487906490781145
$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 29.294891
$ws.Range("N8").Value = 87.88467299999999
$ws.Range("O8").Value = 0.1336376186888105
$ws.Range("P8").Value = 0.1336376186888105
$ws.Range("Q8").Value = 3882.884345645978
$ws.Range("R8").Value = 34945.9591108138
$ws.Range("S8").Value = 0.06722357882159595
$ws.Range("T8").Value = 0.06722357882159596
$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("M9").Value = 67.71760166666667
$ws.Range("N9").Value = 203.152805
$ws.Range("O9").Value = 0.3089145827526977
$ws.Range("P9").Value = 0.3089145827526977
$ws.Range("Q9").Value = 8975.613373546603
$ws.Range("R9").Value = 80780.52036191942
$ws.Range("S9").Value = 0.1553929500283379
$ws.Range("T9").Value = 0.1553929500283379
$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 57.35848733333334
$ws.Range("N10").Value = 172.075462
$ws.Range("O10").Value = 0.261658309594631
$ws.Range("P10").Value = 0.261658309594631
$ws.Range("Q10").Value = 2405.921501635805
$ws.Range("R10").Value = 21653.29351472225
$ws.Range("S10").Value = 0.04165322458938189
$ws.Range("T10").Value = 0.04165322458938189
$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.2957894889638607
$ws.Range("P11").Value = 0.2957894889638607
$ws.Range("Q11").Value = 2719.754218998522
$ws.Range("R11").Value = 24477.7879709867
$ws.Range("S11").Value = 0.04708654593877647
$ws.Range("T11").Value = 0.04708654593877647
$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 29.294891
$ws.Range("N12").Value = 87.88467299999999
$ws.Range("O12").Value = 0.1336376186888105
$ws.Range("P12").Value = 0.1336376186888105
$ws.Range("Q12").Value = 1228.784290202467
$ws.Range("R12").Value = 11059.05861182221
$ws.Range("S12").Value = 0.02127368992583839
$ws.Range("T12").Value = 0.02127368992583839
$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("M13").Value = 67.71760166666667
$ws.Range("N13").Value = 203.152805
$ws.Range("O13").Value = 0.3089145827526977
$ws.Range("P13").Value = 0.3089145827526977
$ws.Range("Q13").Value = 2840.43811933584
$ws.Range("R13").Value = 25563.94307402256
$ws.Range("S13").Value = 0.04917592150720423
$ws.Range("T13").Value = 0.04917592150720423
$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 57.35848733333334
$ws.Range("N14").Value = 172.075462
$ws.Range("O14").Value = 0.261658309594631
$ws.Range("P14").Value = 0.261658309594631
$ws.Range("Q14").Value = 3803.016153965219
$ws.Range("R14").Value = 34227.14538568697
$ws.Range("S14").Value = 0.06584083723033268
$ws.Range("T14").Value = 0.06584083723033268
$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.2957894889638607
$ws.Range("P15").Value = 0.2957894889638607
$ws.Range("Q15").Value = 4299.088404436291
$ws.Range("R15").Value = 38691.79563992662
$ws.Range("S15").Value = 0.07442923416987653
$ws.Range("T15").Value = 0.07442923416987653
$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 29.294891
$ws.Range("N16").Value = 87.88467299999999
$ws.Range("O16").Value = 0.1336376186888105
$ws.Range("P16").Value = 0.1336376186888105
$ws.Range("Q16").Value = 1942.327088477908
$ws.Range("R16").Value = 17480.94379630117
$ws.Range("S16").Value = 0.03362710977369924
$ws.Range("T16").Value = 0.03362710977369924
$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("M17").Value = 67.71760166666667
$ws.Range("N17").Value = 203.152805
$ws.Range("O17").Value = 0.3089145827526977
$ws.Range("P17").Value = 0.3089145827526977
$ws.Range("Q17").Value = 4489.852243653113
$ws.Range("R17").Value = 40408.67019287802
$ws.Range("S17").Value = 0.07773188932010837
$ws.Range("T17").Value = 0.07773188932010837
